# edit.ps1 — reproduce the target diff via PowerPoint COM-interop.
#
# The diff contains two independent changes:
#
#   1. Slide 16 (the "PLENARY - COMPLETE THE MISSING GAPS" slide) has a
#      2-column table whose <a:tableStyleId> changes from the locally
#      defined "Table_0" style ({4DAA1766-...}) to a built-in PowerPoint
#      table style ({A2AB6E76-...}) — i.e. someone picked a different
#      style from the Table Styles gallery on the Table Design ribbon.
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml swap contents: the
#      slide master's theme ("Integral") and the notes master's theme
#      ("Office Theme") trade places. The part that is reachable from the
#      PowerPoint object model is the slide master's theme (also the
#      presentation's primary theme), so we recolour it with the Office
#      Theme's 12 theme colours (the font scheme / format scheme are
#      already identical between the two themes, so only the colour
#      scheme actually needs to change).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style on slide 16
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(16)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{A2AB6E76-48ED-4EAA-BD82-BBB9893E8CAA}")
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours: recolour the slide master's theme from "Integral"
#    to the stock "Office Theme" palette.
# ---------------------------------------------------------------------
function ToCOMRGB($r, $g, $b) {
    # VBA/PowerPoint RGB() packs as 0x00BBGGRR
    return ($b * 65536) + ($g * 256) + $r
}

$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

$colors.Colors(1).RGB  = ToCOMRGB 0x00 0x00 0x00   # dk1
$colors.Colors(2).RGB  = ToCOMRGB 0xFF 0xFF 0xFF   # lt1
$colors.Colors(3).RGB  = ToCOMRGB 0x44 0x54 0x6A   # dk2
$colors.Colors(4).RGB  = ToCOMRGB 0xE7 0xE6 0xE6   # lt2
$colors.Colors(5).RGB  = ToCOMRGB 0x5B 0x9B 0xD5   # accent1
$colors.Colors(6).RGB  = ToCOMRGB 0xED 0x7D 0x31   # accent2
$colors.Colors(7).RGB  = ToCOMRGB 0xA5 0xA5 0xA5   # accent3
$colors.Colors(8).RGB  = ToCOMRGB 0xFF 0xC0 0x00   # accent4
$colors.Colors(9).RGB  = ToCOMRGB 0x44 0x72 0xC4   # accent5
$colors.Colors(10).RGB = ToCOMRGB 0x70 0xAD 0x47   # accent6
$colors.Colors(11).RGB = ToCOMRGB 0x05 0x63 0xC1   # hlink
$colors.Colors(12).RGB = ToCOMRGB 0x95 0x4F 0x72   # folHlink
